{"js": "// Remove the stray \"_GoBack\" bookmark that Word leaves behind after the\n// last edit position, then replace the placeholder \">>> your stuff after\n// this line >>>\" paragraph with the author's own comment.\n\n// 1) Drop the _GoBack bookmark (both its start and end markers).\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// 2) Find the placeholder paragraph (\"...your stuff after this line...\")\n//    and swap its text for the new sentence, keeping the leading/trailing\n//    \">>>\" markers intact.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst target = paragraphs.items.find((p) => p.text.indexOf(\"stuff after this line\") !== -1);\n\nif (target) {\n  target.insertText(\n    \">>> I am Chalani. This is Code/file version management part of my Managing Software Development assignment. >>>\",\n    \"Replace\"\n  );\n  await context.sync();\n}\n", "ps1": "# Remove the stray \"_GoBack\" bookmark that Word leaves behind after the\n# last edit position, then replace the placeholder \">>> your stuff after\n# this line >>>\" paragraph with the author's own comment.\n\n$d = $word.ActiveDocument\n\n# 1) Drop the _GoBack bookmark (removes both its start and end markers).\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n# 2) Find the placeholder text and replace it with the new sentence,\n#    keeping the leading/trailing \">>>\" markers intact.\n$find = $d.Content.Find\n$find.Text = \">>>  your stuff after this line >>>\"\n$find.Replacement.Text = \">>> I am Chalani. This is Code/file version management part of my Managing Software Development assignment. >>>\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n"}
